$d = $word.ActiveDocument

$replacements = @(
    @{old="99×59="; new="55×97="},
    @{old="55×34="; new="31×26="},
    @{old="87×32="; new="68×61="},
    @{old="32×56="; new="52×52="},
    @{old="47×73="; new="48×31="},
    @{old="36×61="; new="74×85="},
    @{old="56×14="; new="72×17="},
    @{old="70×19="; new="87×26="},
    @{old="79×88="; new="52×43="},
    @{old="75×56="; new="24×43="},
    @{old="48×13="; new="48×49="},
    @{old="26×70="; new="90×92="},
    @{old="11×14="; new="91×82="},
    @{old="84×93="; new="82×50="},
    @{old="40×65="; new="19×11="},
    @{old="18×22="; new="70×49="},
    @{old="47×76="; new="87×70="},
    @{old="59×76="; new="14×78="},
    @{old="59×28="; new="24×58="},
    @{old="91×50="; new="38×99="},
    @{old="74×81="; new="21×99="},
    @{old="54×63="; new="47×50="},
    @{old="61×97="; new="25×51="},
    @{old="92×15="; new="97×58="},
    @{old="52×17="; new="29×83="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
